$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 0.018508
$ws.Range("H2").Value = 0.055524
$ws.Range("I2").Value = 0.2347866901774728
$ws.Range("J2").Value = 0.3151800006811757
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2823496666666667
$ws.Range("N2").Value = 0.8470489999999999
$ws.Range("O2").Value = 0.1292103628953711
$ws.Range("P2").Value = 0.1411524606683174
$ws.Range("Q2").Value = 0.005225727630666667
$ws.Range("R2").Value = 0.04703154867599999
$ws.Range("S2").Value = 0.03033687344083433
$ws.Range("T2").Value = 0.0444884326495899

$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("G3").Value = 0.018508
$ws.Range("H3").Value = 0.055524
$ws.Range("I3").Value = 0.2347866901774728
$ws.Range("J3").Value = 0.3151800006811757
$ws.Range("O3").Value = 0.4195646427067258
$ws.Range("P3").Value = 0.4583423527370912
$ws.Range("Q3").Value = 0.01696868964
$ws.Range("R3").Value = 0.15271820676
$ws.Range("S3").Value = 0.09850819377660609
$ws.Range("T3").Value = 0.1444603430478881

$ws.Range("D4").Value = "ECs"
$ws.Range("G4").Value = 0.018508
$ws.Range("H4").Value = 0.055524
$ws.Range("I4").Value = 0.2347866901774728
$ws.Range("J4").Value = 0.3151800006811757
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.19989
$ws.Range("N4").Value = 0.5996699999999999
$ws.Range("O4").Value = 0.09147472969977793
$ws.Range("P4").Value = 0.09992916122794536
$ws.Range("Q4").Value = 0.00369956412
$ws.Range("R4").Value = 0.03329607707999999
$ws.Range("S4").Value = 0.02147704902108983
$ws.Range("T4").Value = 0.03149567310389313

$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 0.018508
$ws.Range("H5").Value = 0.055524
$ws.Range("I5").Value = 0.2347866901774728
$ws.Range("J5").Value = 0.3151800006811757
$ws.Range("M5").Value = 0.55463
$ws.Range("N5").Value = 1.10926
$ws.Range("O5").Value = 0.253812743675961
$ws.Range("P5").Value = 0.1848473683587818
$ws.Range("Q5").Value = 0.01026509204
$ws.Range("R5").Value = 0.06159055223999999
$ws.Range("S5").Value = 0.05959185401254216
$ws.Range("T5").Value = 0.05826019368523438

$ws.Range("D6").Value = "MuSCs"
$ws.Range("G6").Value = 0.018508
$ws.Range("H6").Value = 0.055524
$ws.Range("I6").Value = 0.2347866901774728
$ws.Range("J6").Value = 0.3151800006811757
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.231494
$ws.Range("N6").Value = 0.694482
$ws.Range("O6").Value = 0.1059375210221642
$ws.Range("P6").Value = 0.1157286570078643
$ws.Range("Q6").Value = 0.004284490952
$ws.Range("R6").Value = 0.038560418568
$ws.Range("S6").Value = 0.02487271992640037
$ws.Range("T6").Value = 0.03647535819457021

$ws.Range("D7").Value = "FAPs"
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.060321
$ws.Range("H7").Value = 0.120642
$ws.Range("I7").Value = 0.7652133098225272
$ws.Range("J7").Value = 0.6848199993188243
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2823496666666667
$ws.Range("N7").Value = 0.8470489999999999
$ws.Range("O7").Value = 0.1292103628953711
$ws.Range("P7").Value = 0.1411524606683174
$ws.Range("Q7").Value = 0.017031614243
$ws.Range("R7").Value = 0.102189685458
$ws.Range("S7").Value = 0.09887348945453681
$ws.Range("T7").Value = 0.0966640280187275

$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.060321
$ws.Range("H8").Value = 0.120642
$ws.Range("I8").Value = 0.7652133098225272
$ws.Range("J8").Value = 0.6848199993188243
$ws.Range("O8").Value = 0.4195646427067258
$ws.Range("P8").Value = 0.4583423527370912
$ws.Range("Q8").Value = 0.05530410243
$ws.Range("R8").Value = 0.33182461458
$ws.Range("S8").Value = 0.3210564489301197
$ws.Range("T8").Value = 0.3138820096892032

$ws.Range("D9").Value = "ECs"
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.060321
$ws.Range("H9").Value = 0.120642
$ws.Range("I9").Value = 0.7652133098225272
$ws.Range("J9").Value = 0.6848199993188243
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.19989
$ws.Range("N9").Value = 0.5996699999999999
$ws.Range("O9").Value = 0.09147472969977793
$ws.Range("P9").Value = 0.09992916122794536
$ws.Range("Q9").Value = 0.01205756469
$ws.Range("R9").Value = 0.07234538813999999
$ws.Range("S9").Value = 0.0699976806786881
$ws.Range("T9").Value = 0.06843348812405223

$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.060321
$ws.Range("H10").Value = 0.120642
$ws.Range("I10").Value = 0.7652133098225272
$ws.Range("J10").Value = 0.6848199993188243
$ws.Range("M10").Value = 0.55463
$ws.Range("N10").Value = 1.10926
$ws.Range("O10").Value = 0.253812743675961
$ws.Range("P10").Value = 0.1848473683587818
$ws.Range("Q10").Value = 0.03345583622999999
$ws.Range("R10").Value = 0.13382334492
$ws.Range("S10").Value = 0.1942208896634188
$ws.Range("T10").Value = 0.1265871746735474

$ws.Range("D11").Value = "MuSCs"
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.060321
$ws.Range("H11").Value = 0.120642
$ws.Range("I11").Value = 0.7652133098225272
$ws.Range("J11").Value = 0.6848199993188243
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.231494
$ws.Range("N11").Value = 0.694482
$ws.Range("O11").Value = 0.1059375210221642
$ws.Range("P11").Value = 0.1157286570078643
$ws.Range("Q11").Value = 0.013963949574
$ws.Range("R11").Value = 0.08378369744400001
$ws.Range("S11").Value = 0.08106480109576381
$ws.Range("T11").Value = 0.07925329881329406

